# StagingTemplates/Staging.Group.xlsx: add the header row used by the
# staging import (BusinessKey / Code / Group_ID / Name), bold+underlined,
# on row 2 of Sheet1 (row 1 holds the "for internal use only" banner).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "BusinessKey"
$ws.Range("B2").Value = "Code"
$ws.Range("C2").Value = "Group_ID"
$ws.Range("D2").Value = "Name"

$headerRange = $ws.Range("A2:D2")
$headerRange.Font.Bold = $true
$headerRange.Font.Underline = $true
